# Updates the cryptos list (Price / Volume(1h) columns, plus two name/link
# swaps) to match the latest GitHub Actions scrape.
#
# NOTE: every value is written with a leading apostrophe. That's Excel's
# standard "force text" marker (it is stripped from the stored value but
# keeps the cell typed as Text) so price strings that look numeric, e.g.
# "6.10" or "1.00", don't get silently coerced into numbers and lose their
# trailing zeros - matching the original inlineStr/text cells in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.810.34"
$ws.Range("E2").Value = "'  +2.62%  "
$ws.Range("D3").Value = "'2.230.52"
$ws.Range("E3").Value = "'  +0.72%  "
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("D5").Value = "'231.58"
$ws.Range("E5").Value = "'  +1.02%  "
$ws.Range("D6").Value = "'0.625"
$ws.Range("E6").Value = "'  -1.48%  "
$ws.Range("D7").Value = "'60.42"
$ws.Range("E7").Value = "'  -5.72%  "
$ws.Range("E8").Value = "'  +0.08%  "
$ws.Range("E9").Value = "'  -0.19%  "
$ws.Range("D10").Value = "'58.12"
$ws.Range("E10").Value = "'  -1.90%  "
$ws.Range("D11").Value = "'0.0905"
$ws.Range("E11").Value = "'  +4.26%  "
$ws.Range("E12").Value = "'  -0.19%  "
$ws.Range("D13").Value = "'2.563.22"
$ws.Range("E13").Value = "'  +0.81%  "
$ws.Range("D14").Value = "'15.74"
$ws.Range("E14").Value = "'  -0.76%  "
$ws.Range("D15").Value = "'22.93"
$ws.Range("E15").Value = "'  +2.67%  "
$ws.Range("D16").Value = "'0.803"
$ws.Range("E16").Value = "'  -2.57%  "
$ws.Range("D17").Value = "'5.62"
$ws.Range("E17").Value = "'  +0.11%  "
$ws.Range("D18").Value = "'2.252.53"
$ws.Range("E18").Value = "'  +1.69%  "
$ws.Range("D19").Value = "'41.771.40"
$ws.Range("E19").Value = "'  +2.95%  "
$ws.Range("D20").Value = "'0.0₃0905"
$ws.Range("E20").Value = "'  -0.41%  "
$ws.Range("D21").Value = "'72.33"
$ws.Range("E21").Value = "'  -2.06%  "
$ws.Range("D22").Value = "'6.10"
$ws.Range("E22").Value = "'  -0.30%  "
$ws.Range("D23").Value = "'248.51"
$ws.Range("E23").Value = "'  -0.89%  "
$ws.Range("E24").Value = "'  +0.01%  "
$ws.Range("B25").Value = "'Toncoin"
$ws.Range("C25").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.39"
$ws.Range("E25").Value = "'  +3.50%  "
$ws.Range("B26").Value = "'PancakeSwap"
$ws.Range("C26").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.37"
$ws.Range("E26").Value = "'  +0.19%  "
$ws.Range("E27").Value = "'  +0.32%  "
$ws.Range("D28").Value = "'169.71"
$ws.Range("E28").Value = "'  -1.97%  "
$ws.Range("D29").Value = "'0.142"
$ws.Range("E29").Value = "'  +1.26%  "
$ws.Range("E30").Value = "'  -2.20%  "
$ws.Range("D31").Value = "'1.41"
$ws.Range("E31").Value = "'  -1.94%  "
$ws.Range("E32").Value = "'  -5.62%  "
$ws.Range("E33").Value = "'  -1.40%  "
$ws.Range("D34").Value = "'5.07"
$ws.Range("E34").Value = "'  +6.29%  "
$ws.Range("E35").Value = "'  +0.18%  "
$ws.Range("D36").Value = "'0.0647"
$ws.Range("E36").Value = "'  +2.73%  "
$ws.Range("E37").Value = "'  -7.68%  "
$ws.Range("D38").Value = "'3.64"
$ws.Range("E38").Value = "'  -4.69%  "
$ws.Range("D39").Value = "'2.38"
$ws.Range("E39").Value = "'  -3.92%  "
$ws.Range("B40").Value = "'BinanceUSD"
$ws.Range("C40").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "'  +0.02%  "
$ws.Range("B41").Value = "'TerraClassic"
$ws.Range("C41").Value = "'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D41").Value = "'0.000236"
$ws.Range("E41").Value = "'  +13.83%  "
$ws.Range("D42").Value = "'0.0241"
$ws.Range("E42").Value = "'  +3.70%  "
$ws.Range("E43").Value = "'  -1.02%  "
$ws.Range("E44").Value = "'  -1.30%  "
$ws.Range("D45").Value = "'98.44"
$ws.Range("E45").Value = "'  -3.05%  "
$ws.Range("D46").Value = "'0.0959"
$ws.Range("E46").Value = "'  +2.15%  "
$ws.Range("D47").Value = "'4.42"
$ws.Range("E47").Value = "'  -8.57%  "
$ws.Range("D48").Value = "'1.469.33"
$ws.Range("E48").Value = "'  -3.42%  "
$ws.Range("D49").Value = "'16.58"
$ws.Range("E49").Value = "'  -4.09%  "
$ws.Range("E50").Value = "'  -1.45%  "
$ws.Range("D51").Value = "'2.27"
$ws.Range("E51").Value = "'  +10.17%  "
